# Auto-generated edit script applying the diff changes to the workbook.
# Numeric cells (columns F, G) are assigned directly as numbers.
# Text cells are protected against Excel's automatic date/number
# auto-detection (e.g. "2024-05-01" -> date serial) by temporarily
# forcing a Text number format, then resetting the style back to
# "Normal" afterwards so no stray style index is left behind.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 351
$ws.Range("F3").Value = 1253
$ws.Range("F5").Value = 78
$ws.Range("F6").Value = 228
$ws.Range("F7").Value = 715
$ws.Range("F8").Value = 2036
$ws.Range("F10").Value = 727
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "上海·坏孩纸物语第40届动漫节之曹沫篇"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "曹杨路1888号 复悦荟"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2024.05.01 11:00-05.01 17:00"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = 529
$ws.Range("G12").Value = 66.90000000000001
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84724"
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202404/Y4wTU9111713328435995.png"
$ws.Range("I12").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "上海·次元裂缝-X Anikura二次元派对"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024.05.01 17:00-05.01 22:00"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = 101
$ws.Range("G13").Value = 60
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=84409"
$ws.Range("H13").Style = "Normal"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202404/z38QIjBW1713260568891.jpeg"
$ws.Range("I13").Style = "Normal"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "上海·百梦动漫游戏嘉年华"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "吴中路1588号上海爱琴海购物中心F4 百忍潮玩对战中心(爱琴海购物中心店)"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2024.05.01 10:00-05.02 19:00"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = 129
$ws.Range("G14").Value = 66
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=84152"
$ws.Range("H14").Style = "Normal"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202404/hmAEptJH1713249167991.jpeg"
$ws.Range("I14").Style = "Normal"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "上海·第一届妖妖动漫游戏展"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2024.05.01 10:00-05.04 17:00"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = 1049
$ws.Range("G15").Value = 68
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=84642"
$ws.Range("H15").Style = "Normal"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202404/fGytR92V1714112934007.jpeg"
$ws.Range("I15").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "上海·第五十七届燃梦星辰动漫嘉年华"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "云锦路500号(近11号线地铁站5号口) 绿地滨江CLUB"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2024.05.01 10:30-05.01 16:30"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = 779
$ws.Range("G16").Value = 58.8
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=83807"
$ws.Range("H16").Style = "Normal"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202404/RGLpPX211712156496032.jpeg"
$ws.Range("I16").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "上海·第十一届ACBC·妖妖动漫游戏展"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "莫干山路600号 大洋晶典天安千树"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2024.05.01 11:00-05.01 18:00"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 48.8
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=84765"
$ws.Range("H17").Style = "Normal"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202404/7LXcHPbn1713936506534.jpeg"
$ws.Range("I17").Style = "Normal"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "2024-05-01"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "上海·魔都劳动节漫展-CF01"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "澳门路168号 月星家居（澳门路）"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2024.05.01 10:00-05.05 16:00"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = 623
$ws.Range("G18").Value = 49
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=82992"
$ws.Range("H18").Style = "Normal"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202403/I7O9LMtb1710752670542.jpeg"
$ws.Range("I18").Style = "Normal"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "上海·2024GAF插画艺术节"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "博成路850号 上海世博展览馆"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2024.05.02 10:30-05.04 19:00"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Value = 1178
$ws.Range("G19").Value = 128
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=83699"
$ws.Range("H19").Style = "Normal"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202403/APlNld8y1711825700811.jpeg"
$ws.Range("I19").Style = "Normal"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "上海·坏孩子物语第37届动漫节之“要离”篇"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "沪南路2229号 复地活力城"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2024.05.02 13:00-05.02 17:40"
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 66.90000000000001
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=84811"
$ws.Range("H20").Style = "Normal"
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202404/fFfuNGUu1713967452154.png"
$ws.Range("I20").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "上海·女团驾到·次元女团偶像专区"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2024.05.02 13:00-05.04 18:00"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 78
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=84796"
$ws.Range("H21").Style = "Normal"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202404/AOS8NlZ31713944402838.jpeg"
$ws.Range("I21").Style = "Normal"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "2024-05-02"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "周家嘴路3608号 宝龙旭辉广场"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2024.05.02 10:20-05.03 16:30"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Value = 722
$ws.Range("G22").Value = 58
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=82761"
$ws.Range("H22").Style = "Normal"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202403/azEA4EM01710236719279.jpeg"
$ws.Range("I22").Style = "Normal"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "上海·DizzyMart2024电则市场 中国同人音乐展会"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "顾村镇蕰川路6号 智慧湾科创园"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2024.05.03 09:00-05.04 20:00"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = 686
$ws.Range("G23").Value = 138
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=84202"
$ws.Range("H23").Style = "Normal"
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202404/etRgMvxv1712656961255.jpeg"
$ws.Range("I23").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "上海·HD动漫主题嘉年华·大唐夜市之剑侠奇缘"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "天等路400号，与华东理工大学仅一墙之隔 品域凌云里"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2024.05.03 10:00-05.04 17:30"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 75
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=84247"
$ws.Range("H24").Style = "Normal"
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = "//i2.hdslb.com/bfs/openplatform/202404/WqEbY0fn1714099426677.jpeg"
$ws.Range("I24").Style = "Normal"
$ws.Range("F26").Value = 613
$ws.Range("F27").Value = 1178
$ws.Range("F30").Value = 4703
$ws.Range("F32").Value = 1363
$ws.Range("F33").Value = 5722
$ws.Range("F34").Value = 940
$ws.Range("F35").Value = 559
$ws.Range("F36").Value = 51
$ws.Range("F38").Value = 1030
$ws.Range("F41").Value = 627
$ws.Range("F47").Value = 84
$ws.Range("F49").Value = 9

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 2063
$ws.Range("F10").Value = 458
$ws.Range("F12").Value = 89
$ws.Range("F38").Value = 54
$ws.Range("F43").Value = 461
$ws.Range("F45").Value = 76

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 645

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 351
$ws.Range("F4").Value = 645
$ws.Range("F6").Value = 1253
$ws.Range("F8").Value = 78
$ws.Range("F9").Value = 228
$ws.Range("F10").Value = 715
$ws.Range("F14").Value = 2063
$ws.Range("F15").Value = 2036
$ws.Range("F17").Value = 727
$ws.Range("F21").Value = 129
$ws.Range("F22").Value = 1049
$ws.Range("F23").Value = 779
$ws.Range("F24").Value = 18
$ws.Range("F25").Value = 1178
$ws.Range("F26").Value = 89
$ws.Range("F28").Value = 722
$ws.Range("F30").Value = 686
$ws.Range("F32").Value = 613
$ws.Range("F38").Value = 4703
$ws.Range("F39").Value = 1363
$ws.Range("F40").Value = 5722
$ws.Range("F41").Value = 940
$ws.Range("F43").Value = 560
$ws.Range("F44").Value = 51
$ws.Range("F45").Value = 1030
$ws.Range("F46").Value = 627
$ws.Range("F47").Value = 54
$ws.Range("F51").Value = 461
$ws.Range("F52").Value = 84

